# Break out stock.yaml completed: refresh the "backup" (R) column for the
# already-fetched rows and append the newly-fetched trading days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) R1317:R1354 were left as blank placeholders while rows were still
#        being processed. Now that processing finished, mirror the
#        "detect_structure" (Q) value into "backup" (R) for each of them.
for ($r = 1317; $r -le 1354; $r++) {
    $detectStructure = $ws.Cells.Item($r, 17).Value2
    $ws.Cells.Item($r, 18).Value = $detectStructure
}

# --- 2) Append the 13 newly scraped trading days (2024-10-08 .. 2024-10-24).
#        Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close,
#        G Volume, H Year, I Month, J Day, K Hour, L Minute, M Second,
#        N Week, O isPivot, P two_line_structure, Q detect_structure.
#        R (backup) is intentionally left blank for these brand-new rows,
#        matching the same "not yet processed" state rows 1317-1354 had
#        before step 1 above.
$newRows = @(
    @(1355, 45573, 2175,              2284,              2172.60009765625,  2276.89990234375, 2276.89990234375, 114753, 2024, 10, 8,  0, 0, 0, 41, 0, 0, 1),
    @(1356, 45574, 2270,              2312.35009765625,  2268.050048828125, 2288.550048828125, 2288.550048828125, 105611, 2024, 10, 9,  0, 0, 0, 41, 0, 0, 0),
    @(1357, 45575, 2308,              2318.300048828125, 2226,              2259.550048828125, 2259.550048828125, 96249,  2024, 10, 10, 0, 0, 0, 41, 1, 0, 0),
    @(1358, 45576, 2260,              2295.89990234375,  2236.300048828125, 2274.050048828125, 2274.050048828125, 132443, 2024, 10, 11, 0, 0, 0, 41, 0, 0, 0),
    @(1359, 45579, 2286.75,           2303.300048828125, 2212.550048828125, 2288.449951171875, 2288.449951171875, 127726, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(1360, 45580, 2283,              2289.050048828125, 2235.10009765625,  2259.550048828125, 2259.550048828125, 66975,  2024, 10, 15, 0, 0, 0, 42, 0, 0, 0),
    @(1361, 45581, 2249.800048828125, 2269.25,            2218.60009765625,  2228.699951171875, 2228.699951171875, 38676,  2024, 10, 16, 0, 0, 0, 42, 0, 0, 0),
    @(1362, 45582, 2237,              2261.5,             2178.550048828125, 2187.800048828125, 2187.800048828125, 187850, 2024, 10, 17, 0, 0, 0, 42, 0, 0, 0),
    @(1363, 45583, 2189.60009765625,  2218,               2175.75,           2205.64990234375,  2205.64990234375,  71188,  2024, 10, 18, 0, 0, 0, 42, 0, 0, 0),
    @(1364, 45586, 2200.85009765625,  2217.89990234375,   2170.5,            2176.5,             2176.5,             27455,  2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(1365, 45587, 2189.85009765625,  2190.5,             2127,              2133.300048828125,  2133.300048828125,  62564,  2024, 10, 22, 0, 0, 0, 43, 0, 0, 0),
    @(1366, 45588, 2129.050048828125, 2173.699951171875,  2116.449951171875, 2121.449951171875,  2121.449951171875,  87491,  2024, 10, 23, 0, 0, 0, 43, 0, 0, 0),
    @(1367, 45589, 2134.949951171875, 2165.050048828125,  2109.5,            2135.35009765625,   2135.35009765625,   79476,  2024, 10, 24, 0, 0, 0, 43, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Datetime (column A) keeps the same date/time number format as the
    # existing rows above it.
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # B..Q = Open, High, Low, Close, Adj Close, Volume, Year, Month, Day,
    #        Hour, Minute, Second, Week, isPivot, two_line_structure,
    #        detect_structure
    for ($c = 2; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
